$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to stay text while assigning new values, to avoid Excel
# auto-converting plain-decimal-looking strings (e.g. "19.68") into numbers.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = "26.638.14"
$ws.Range("E2").Value = "  -1.58%  "
$ws.Range("D3").Value = "1.591.69"
$ws.Range("E3").Value = "  -1.99%  "
$ws.Range("E4").Value = "  -0.25%  "
$ws.Range("D5").Value = "211.22"
$ws.Range("E5").Value = "  -1.87%  "
$ws.Range("E6").Value = "  -1.10%  "
$ws.Range("E7").Value = "  -0.38%  "
$ws.Range("E8").Value = "  -2.34%  "
$ws.Range("D9").Value = "0.0616"
$ws.Range("E9").Value = "  -1.25%  "
$ws.Range("D10").Value = "19.68"
$ws.Range("E10").Value = "  -2.20%  "
$ws.Range("D11").Value = "0.0835"
$ws.Range("E11").Value = "  -1.50%  "
$ws.Range("E12").Value = "  -2.02%  "
$ws.Range("D13").Value = "1.591.26"
$ws.Range("E13").Value = "  -2.04%  "
$ws.Range("E14").Value = "  -2.01%  "
$ws.Range("D15").Value = "0.527"
$ws.Range("E15").Value = "  -2.68%  "
$ws.Range("D16").Value = "64.78"
$ws.Range("E16").Value = "  +0.53%  "
$ws.Range("D17").Value = "26.634.54"
$ws.Range("E17").Value = "  -1.46%  "
$ws.Range("E18").Value = "  -1.52%  "
$ws.Range("D19").Value = "208.47"
$ws.Range("E19").Value = "  -3.49%  "
$ws.Range("E20").Value = "  -0.21%  "
$ws.Range("D21").Value = "6.75"
$ws.Range("E22").Value = "  -2.64%  "
$ws.Range("E23").Value = "  -2.30%  "
$ws.Range("D24").Value = "8.90"
$ws.Range("E24").Value = "  -0.82%  "
$ws.Range("D25").Value = "146.69"
$ws.Range("E25").Value = "  -0.64%  "
$ws.Range("E26").Value = "  -0.34%  "
$ws.Range("D27").Value = "7.27"
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("E28").Value = "  -3.36%  "
$ws.Range("E29").Value = "  -1.62%  "
$ws.Range("E30").Value = "  +0.59%  "
$ws.Range("E31").Value = "  -1.76%  "
$ws.Range("E32").Value = "  -3.46%  "
$ws.Range("E33").Value = "  +22.28%  "
$ws.Range("E34").Value = "  -2.47%  "
$ws.Range("D35").Value = "1.319.27"
$ws.Range("E35").Value = "  -1.30%  "
$ws.Range("E36").Value = "  -4.17%  "
$ws.Range("D37").Value = "2.40"
$ws.Range("E37").Value = "  -2.45%  "
$ws.Range("D38").Value = "0.0172"
$ws.Range("E38").Value = "  -2.04%  "
$ws.Range("D39").Value = "0.829"
$ws.Range("E39").Value = "  -1.99%  "
$ws.Range("E40").Value = "  -0.31%  "
$ws.Range("D41").Value = "5.40"
$ws.Range("E41").Value = "  +3.77%  "
$ws.Range("E42").Value = "  -1.58%  "
$ws.Range("E43").Value = "  -3.20%  "
$ws.Range("D44").Value = "63.31"
$ws.Range("E44").Value = "  -1.48%  "
$ws.Range("D45").Value = "1.726.52"
$ws.Range("E45").Value = "  -2.07%  "
$ws.Range("D46").Value = "90.02"
$ws.Range("E46").Value = "  -0.40%  "
$ws.Range("E47").Value = "  +0.15%  "
$ws.Range("E48").Value = "  +2.13%  "
$ws.Range("E49").Value = "  -0.31%  "
$ws.Range("D50").Value = "0.0978"
$ws.Range("E50").Value = "  -0.60%  "
$ws.Range("D51").Value = "7.55"
$ws.Range("E51").Value = "  +0.07%  "

# Restore original (default) style on column D now that text values are set,
# matching the workbook's original formatting (no explicit number format).
$dRange.Style = "Normal"
